$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for numeric-looking price cells to preserve exact text representation
$textCells = @("D5", "D10", "D11", "D16", "D20", "D21", "D24", "D25", "D26", "D29", "D40", "D43", "D44", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply all cell value updates
$ws.Range('D2').Value = '26.648.93'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.598.34'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '210.76'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D10').Value = '19.58'
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('D11').Value = '0.0846'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '1.822.47'
$ws.Range('D13').Value = '1.580.34'
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('D16').Value = '64.62'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('D17').Value = '26.619.95'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('E18').Value = '  -2.38%  '
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').Value = '208.41'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').Value = '7.09'
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('E23').Value = '  -2.97%  '
$ws.Range('D24').Value = '8.95'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').Value = '143.76'
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +0.45%  '
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('D29').Value = '15.28'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').Value = '  -1.91%  '
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('E34').Value = '  +19.27%  '
$ws.Range('D35').Value = '1.279.91'
$ws.Range('E35').Value = '  -0.75%  '
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('E38').Value = '  -3.48%  '
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('D40').Value = '0.824'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('D43').Value = '0.773'
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('D44').Value = '62.57'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('D45').Value = '1.734.19'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0104'
$ws.Range('E48').Value = '  -2.90%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.103'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0513'
$ws.Range('E50').Value = '  +0.62%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.14%  '
